$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: fill in the Iteration value that was missing (C3)
$ws.Range("C3").Value = 123604301

# Insert a new row at 4, pushing the old row 4 (seed 2, ...) and everything
# below it down by one. The new row 4 holds the one-off
# "vector<int> tabu_tenure_table_i" benchmark result instead.
$ws.Rows("4").Insert()
$ws.Range("A4").Value = "vector<int> tabu_tenure_table_i"
$ws.Range("D4").Value = 843.386

# Insert another blank row at 5 to leave a gap between the new result row
# and the resumed "vector" benchmark series (now starting at row 6).
$ws.Rows("5").Insert()

# The row-insert operations copy the neighbouring cell formatting into the
# new row, which would otherwise leave stray empty-but-styled cells behind.
# Clear those completely so the row matches the saved state exactly.
$ws.Range("C4").Clear()
$ws.Range("C5").Clear()
$ws.Range("D5").Clear()

# Update the active selection to match the saved workbook state.
$ws.Range("F8").Select() | Out-Null
